# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals.
# Only column G ("K") values for rows 2-17 change per the new calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 5
    6  = 8
    7  = 3
    8  = 4
    9  = 6
    10 = 3
    11 = 5
    12 = 3
    13 = 1
    14 = 3
    15 = 3
    16 = 3
    17 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
